# Update trading results - Mon Oct 20 01:33:32 UTC 2025
# Appends two new trading-log rows (90 and 91) to the end of the sheet,
# mirroring the TRADING_ATTEMPT -> POSITION_OPENED pattern used throughout
# the rest of the log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 90: TRADING_ATTEMPT for TRX
$ws.Cells.Item(90, 1).Value = "2025-10-20T01:33:30.469956"
$ws.Cells.Item(90, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(90, 3).Value = "TRX"
$ws.Cells.Item(90, 4).Value = "UNKNOWN"
$ws.Cells.Item(90, 5).Value = 0.3208392693959338
$ws.Cells.Item(90, 11).Value = "ATTEMPT"
$ws.Cells.Item(90, 12).Value = "Attempting trade 1/1"

# Row 91: POSITION_OPENED for TRX
$ws.Cells.Item(91, 1).Value = "2025-10-20T01:33:32.061100"
$ws.Cells.Item(91, 2).Value = "POSITION_OPENED"
$ws.Cells.Item(91, 3).Value = "TRX"
$ws.Cells.Item(91, 4).Value = "UNKNOWN"
$ws.Cells.Item(91, 5).Value = 0.3208392693959338
$ws.Cells.Item(91, 6).Value = 900
$ws.Cells.Item(91, 7).Value = 10
$ws.Cells.Item(91, 8).Value = 0.1367964106814732
$ws.Cells.Item(91, 11).Value = "SUCCESS"
